$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.328.54'
Set-TextValue 'E2' '  -4.52%  '
Set-TextValue 'D3' '1.760.97'
Set-TextValue 'E3' '  -4.14%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'E5' '  +0.01%  '
Set-TextValue 'D6' '301.77'
Set-TextValue 'E6' '  -3.35%  '
Set-TextValue 'D7' '0.4267'
Set-TextValue 'E7' '  -0.33%  '
Set-TextValue 'D8' '0.3616'
Set-TextValue 'E8' '  -1.18%  '
Set-TextValue 'D9' '0.07039'
Set-TextValue 'E9' '  -3.30%  '
Set-TextValue 'D10' '0.8300'
Set-TextValue 'E10' '  -4.06%  '
Set-TextValue 'D11' '20.13'
Set-TextValue 'E11' '  -2.77%  '
Set-TextValue 'D12' '1.763.72'
Set-TextValue 'E12' '  -0.69%  '
Set-TextValue 'E13' '  -4.19%  '
Set-TextValue 'D14' '6.378'
Set-TextValue 'E14' '  -2.26%  '
Set-TextValue 'D15' '0.06787'
Set-TextValue 'E15' '  -2.52%  '
Set-TextValue 'D16' '1.005'
Set-TextValue 'E16' '  +0.21%  '
Set-TextValue 'D17' '79.09'
Set-TextValue 'E17' '  -2.02%  '
Set-TextValue 'D18' '0.000008617'
Set-TextValue 'E18' '  -3.41%  '
Set-TextValue 'E19' '  -0.08%  '
Set-TextValue 'D20' '14.90'
Set-TextValue 'E20' '  -3.46%  '
Set-TextValue 'D21' '25.850.98'
Set-TextValue 'E21' '  -5.71%  '
Set-TextValue 'D22' '4.991'
Set-TextValue 'E22' '  -3.34%  '
Set-TextValue 'E23' '  +1.69%  '
Set-TextValue 'D24' '1.939.01'
Set-TextValue 'E24' '  -4.56%  '
Set-TextValue 'D25' '1.904'
Set-TextValue 'E25' '  -4.42%  '
Set-TextValue 'D26' '151.87'
Set-TextValue 'E26' '  -1.93%  '
Set-TextValue 'D27' '18.12'
Set-TextValue 'E27' '  -3.97%  '
Set-TextValue 'B28' 'BitcoinCash'
Set-TextValue 'C28' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D28' '114.49'
Set-TextValue 'E28' '  +0.14%  '
Set-TextValue 'B29' 'InternetComputer(DFINITY)'
Set-TextValue 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D29' '5.015'
Set-TextValue 'E29' '  -2.95%  '
Set-TextValue 'D30' '1.665'
Set-TextValue 'E30' '  -8.69%  '
Set-TextValue 'D31' '0.08881'
Set-TextValue 'E31' '  +0.21%  '
Set-TextValue 'D32' '0.7225'
Set-TextValue 'E32' '  -3.77%  '
Set-TextValue 'D33' '1.110'
Set-TextValue 'E33' '  -2.31%  '
Set-TextValue 'D34' '4.296'
Set-TextValue 'E34' '  -5.53%  '
Set-TextValue 'D35' '0.9997'
Set-TextValue 'E35' '  -0.08%  '
Set-TextValue 'D36' '2.716'
Set-TextValue 'E36' '  -9.40%  '
Set-TextValue 'D37' '1.064'
Set-TextValue 'E37' '  -3.06%  '
Set-TextValue 'D38' '0.05088'
Set-TextValue 'E38' '  -4.52%  '
Set-TextValue 'D39' '0.01877'
Set-TextValue 'E39' '  -3.15%  '
Set-TextValue 'D40' '0.1600'
Set-TextValue 'E40' '  -3.08%  '
Set-TextValue 'D41' '0.4881'
Set-TextValue 'E41' '  -3.94%  '
Set-TextValue 'D42' '6.151'
Set-TextValue 'E42' '  -5.23%  '
Set-TextValue 'D43' '2.477'
Set-TextValue 'E43' '  -11.46%  '
Set-TextValue 'D44' '7.955'
Set-TextValue 'E44' '  -4.73%  '
Set-TextValue 'D45' '104.49'
Set-TextValue 'E45' '  -1.00%  '
Set-TextValue 'E47' '  -3.63%  '
Set-TextValue 'D48' '0.06183'
Set-TextValue 'E48' '  -4.48%  '
Set-TextValue 'D49' '0.4453'
Set-TextValue 'E49' '  -4.87%  '
Set-TextValue 'D50' '1.565'
Set-TextValue 'E50' '  -3.48%  '
Set-TextValue 'D51' '1.706'
Set-TextValue 'E51' '  -1.97%  '
